$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.282.91"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.668.35"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'219.82"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D6").Value = "'0.5290"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'0.2652"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Value = "'20.96"
$ws.Range("E10").Value = "  +2.47%  "

$ws.Range("D11").Value = "'0.07837"
$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").Value = "'4.526"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "1.669.44"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").Value = "1.896.66"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "'0.5595"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").Value = "0.0₅8097"
$ws.Range("E16").Value = "  -0.97%  "

$ws.Range("D17").Value = "'65.73"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "26.304.74"
$ws.Range("E18").Value = "  +0.27%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "'4.723"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").Value = "'200.46"
$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("E22").Value = "  +0.94%  "

$ws.Range("D23").Value = "'6.061"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").Value = "'1.009"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'146.42"

$ws.Range("D27").Value = "'7.234"
$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").Value = "'16.20"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").Value = "'1.528"
$ws.Range("E29").Value = "  +3.22%  "

$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("D31").Value = "'1.284"
$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("D32").Value = "'3.511"
$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").Value = "'3.335"
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("E34").Value = "  -0.85%  "

$ws.Range("D35").Value = "'0.9642"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").Value = "'2.814"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").Value = "'0.5807"
$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").Value = "'5.986"
$ws.Range("E40").Value = "  +1.60%  "

$ws.Range("D41").Value = "1.078.72"
$ws.Range("E41").Value = "  +3.44%  "

$ws.Range("D42").Value = "'0.8573"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").Value = "'1.007"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").Value = "'102.91"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").Value = "1.807.16"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "'58.56"
$ws.Range("E46").Value = "  +2.55%  "

$ws.Range("D47").Value = "'1.014"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").Value = "'0.4413"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("D49").Value = "'8.057"
$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("E50").Value = "  -4.76%  "

$ws.Range("D51").Value = "'0.05144"
$ws.Range("E51").Value = "  -0.33%  "
